$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the 19th May row (row 4) training plan entries
$ws.Range("D4").Value2 = "10k easy to HR with 2 x 500m at HMP "
$ws.Range("E4").Value2 = "Rest day or 30min very slow jog"
$ws.Range("F4").Value2 = "5k easy with 6 to 8 100m efforts to 4:05 pace."
$ws.Range("H4").Value2 = "10km easy"

$ws.Range("H16").Select()
